$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 182; $r -le 381; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $d = [DateTime]::FromOADate($cell.Value2)
    $d2 = $d.AddYears(1)
    $cell.Value2 = $d2.ToOADate()
}
